# repull data, push all data, mean calculation
# Update the dSF (column F) values for each row with the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -3
    3  = 5
    4  = 6
    6  = -3
    7  = 3
    8  = -1
    9  = -1
    10 = 2
    11 = 1
    12 = 7
    13 = 1
    14 = 4
    15 = -3
    16 = 6
    17 = 2
    18 = 3
    19 = -3
    20 = 4
    22 = -1
    23 = 3
    24 = 2
    25 = 1
    26 = -4
    27 = -3
    28 = -1
    29 = 2
    30 = 2
    31 = 2
    33 = 6
    34 = 2
    35 = -3
    36 = 2
    38 = -2
    39 = 1
    40 = 4
    41 = -3
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
